# Revise build script optimization
# Appends additional sensor-reading rows to each worksheet, matching the
# pattern of the existing rows in that sheet (same B/D/F/G/H/I constants,
# same C/E pair, new timestamps in column A).

$wb = $excel.ActiveWorkbook

function Add-SensorRows($ws, $startRow, $times, $bVal, $cVal, $dVal, $eVal, $fVal, $gVal, $hVal, $iVal) {
    $templateFormat = $ws.Cells.Item($startRow - 1, 1).NumberFormat()

    for ($k = 0; $k -lt $times.Length; $k++) {
        $r = $startRow + $k

        $ws.Cells.Item($r, 1).NumberFormat = $templateFormat
        $ws.Cells.Item($r, 1).Value = $times[$k]
        $ws.Cells.Item($r, 2).Value = $bVal
        $ws.Cells.Item($r, 3).Value = $cVal
        $ws.Cells.Item($r, 4).Value = $dVal
        $ws.Cells.Item($r, 5).Value = $eVal
        $ws.Cells.Item($r, 6).Value = $fVal
        $ws.Cells.Item($r, 7).Value = $gVal
        $ws.Cells.Item($r, 8).Value = $hVal
        $ws.Cells.Item($r, 9).Value = $iVal
    }
}

$timesA = @(45729.08020857639, 45729.08023003472, 45729.08025329861)
$timesB = @(45729.06475774306, 45729.06477990741, 45729.06480305555, 45729.14823328704, 45729.14825528935, 45729.14827864584, 45729.23170871528, 45729.23173081018, 45729.23175395833)

# the G column value is stored in scientific notation in the source data;
# this mini-shell's parser does not accept an `e+NN` exponent literal, so
# build the double by casting the literal string instead.
$gConst = [double]"5.686312626471138e+23"

# Sheet 1: ROW50-FE-LIFTER (rows 26-28)
$ws1 = $wb.Worksheets.Item(1)
Add-SensorRows $ws1 26 $timesA "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c," "0x01,0x90," "0x14" 400 $gConst 400 20

# Sheet 2: ROW50-MID-LIFTER (rows 59-67)
$ws2 = $wb.Worksheets.Item(2)
Add-SensorRows $ws2 59 $timesB "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x90," "0x19" 400 $gConst 400 25

# Sheet 3: ROW11-FE-LIFTER (rows 26-28)
$ws3 = $wb.Worksheets.Item(3)
Add-SensorRows $ws3 26 $timesA "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c," "0x01,0x90," "0x14" 400 $gConst 400 20

# Sheet 4: ROW11-MID-LIFTER (rows 59-67)
$ws4 = $wb.Worksheets.Item(4)
Add-SensorRows $ws4 59 $timesB "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x90," "0x19" 400 $gConst 400 25

Write-Output "Added rows to all 4 sheets"
